$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.591.89'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.895.33'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.49'
$ws.Range('E5').Value = '  -0.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.693'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.08'
$ws.Range('E8').Value = '  -2.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '57.00'
$ws.Range('E9').Value = '  +9.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.358'
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0753'
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.61'
$ws.Range('E13').Value = '  +10.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.795'
$ws.Range('E14').Value = '  +8.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.171.09'
$ws.Range('E15').Value = '  -0.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.04'
$ws.Range('E16').Value = '  +1.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.898.34'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.574.72'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.64'
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0832'
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '246.39'
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.99'
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.20'
$ws.Range('E23').Value = '  +4.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.68'
$ws.Range('E24').Value = '  +4.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('E26').Value = '  -3.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.91'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.68'
$ws.Range('E28').Value = '  +1.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.39'
$ws.Range('E29').Value = '  -0.28%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.39'
$ws.Range('E31').Value = '  +2.78%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0609'
$ws.Range('E32').Value = '  +4.59%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E34').Value = '  +17.36%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.48'
$ws.Range('E36').Value = '  -17.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.856'
$ws.Range('E38').Value = '  -3.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0737'
$ws.Range('E39').Value = '  +7.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0229'
$ws.Range('E40').Value = '  +6.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.21'
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.09'
$ws.Range('E42').Value = '  -0.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.09'
$ws.Range('E43').Value = '  -0.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.22'
$ws.Range('E44').Value = '  +16.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.317.56'
$ws.Range('E45').Value = '  +1.33%  '
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0810'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.39'
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '42.66'
$ws.Range('E51').Value = '  -2.34%  '
